{"js": "// \"clockwise\" -> \"anti-clockwise\" in the \"All counters are redistributed...\"\n// sentence (Game flow / Counters redistribution paragraph).\n//\n// The source run is split into three runs (matching the canonical OOXML in\n// the diff):\n//   1) \"All counters are redistributed to other pits in \"\n//   2) \"anti-\"\n//   3) \"clockwise direction.\"\n//\n// Plain text search + insertText(before) already performs the split at the\n// insertion point; toggling formatting on the freshly-inserted \"anti-\" text\n// (and resetting it right back) forces the run boundary to persist instead\n// of being re-merged with its neighbours on save, which reproduces the\n// run-per-insertion structure Word itself produces when text is typed in\n// a separate edit.\n\nconst searchResults = context.document.body.search(\"clockwise direction.\", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find \"clockwise direction.\" in the document body.');\n}\n\n// Only one occurrence is expected/handled; insert \"anti-\" right before it.\nconst target = searchResults.items[0];\nconst insertedRange = target.insertText(\"anti-\", Word.InsertLocation.before);\n\n// Force the inserted text to materialize as its own run (distinct from the\n// runs before/after it) by nudging a character property and restoring it.\ninsertedRange.font.bold = true;\ninsertedRange.font.bold = false;\n\nawait context.sync();\n", "ps1": "# \"clockwise\" -> \"anti-clockwise\" in the \"All counters are redistributed...\"\n# sentence (Game flow / Counters redistribution paragraph).\n#\n# The original single run is split into three runs (matching the canonical\n# OOXML in the diff):\n#   1) \"All counters are redistributed to other pits in \"\n#   2) \"anti-\"\n#   3) \"clockwise direction.\"\n#\n# Find.Execute locates the unique \"clockwise direction.\" text, Collapse(1)\n# moves to its start, and InsertBefore(\"anti-\") splits the run there (the\n# range then refers to the freshly inserted \"anti-\" text). Toggling the\n# Font.Bold property on that inserted text and immediately resetting it\n# forces the new run boundary to persist on save instead of being re-merged\n# with its neighbouring runs.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"clockwise direction.\")\nif (-not $found) {\n    throw \"Could not find 'clockwise direction.' in the document.\"\n}\n\n$rng.Collapse(1)          # wdCollapseStart\n$rng.InsertBefore(\"anti-\")\n\n$rng.Font.Bold = 1\n$rng.Font.Bold = 0\n"}
